# Week 15 simulations update
$wb = $excel.ActiveWorkbook

$off = $wb.Worksheets.Item("OFF")
$off.Range("B2").Value = 370
$off.Range("C2").Value = 262
$off.Range("D2").Value = 103
$off.Range("E2").Value = 55
$off.Range("F2").Value = 7
$off.Range("G2").Value = 7

$def = $wb.Worksheets.Item("DEF")
$def.Range("B2").Value = 354
$def.Range("C2").Value = 250
$def.Range("D2").Value = 101
$def.Range("E2").Value = 54
$def.Range("F2").Value = 7
$def.Range("G2").Value = 3
